# Refresh the cryptos price/volume table (GitHub Actions scrape update).
# Values that look numeric are quote-prefixed ('...) so Excel stores them
# as literal text (e.g. "513.10"), matching the source data's formatting
# instead of collapsing to a float (513.1) that would drop trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.847.48"
$ws.Range("E2").Value = "  -2.05%  "
$ws.Range("D3").Value = "2.624.27"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'513.10"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").Value = "'144.08"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("D7").Value = "'0.995"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("D8").Value = "'0.568"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "2.647.99"
$ws.Range("E9").Value = "  +0.85%  "
$ws.Range("D10").Value = "'6.33"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.105"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "'0.337"
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "'0.127"
$ws.Range("E13").Value = "  -1.66%  "
$ws.Range("D14").Value = "3.082.69"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("D15").Value = "58.808.90"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "'21.05"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "'0.0000137"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "2.638.42"
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("D19").Value = "'4.54"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "'343.72"
$ws.Range("E20").Value = "  +0.98%  "
$ws.Range("D21").Value = "'10.35"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").Value = "'6.09"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "'61.01"
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'0.420"
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("D26").Value = "2.739.09"
$ws.Range("E26").Value = "  -0.25%  "
$ws.Range("D27").Value = "'0.992"
$ws.Range("E27").Value = "  -0.74%  "
$ws.Range("E28").Value = "  +0.46%  "
$ws.Range("D29").Value = "0.0₃0805"
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("D30").Value = "'7.13"
$ws.Range("E30").Value = "  +2.14%  "
$ws.Range("D31").Value = "'0.997"
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("D32").Value = "'6.46"
$ws.Range("E32").Value = "  +8.83%  "
$ws.Range("D33").Value = "'1.58"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'18.85"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'150.41"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +12.45%  "
$ws.Range("D37").Value = "'4.02"
$ws.Range("E37").Value = "  +2.70%  "
$ws.Range("D38").Value = "'1.15"
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'0.853"
$ws.Range("E39").Value = "  +0.24%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").Value = "'36.42"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("D41").Value = "'3.69"
$ws.Range("E41").Value = "  +2.24%  "
$ws.Range("D42").Value = "'1.41"
$ws.Range("E42").Value = "  -0.46%  "
$ws.Range("D43").Value = "'280.22"
$ws.Range("E43").Value = "  -2.46%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.612"
$ws.Range("E44").Value = "  -2.30%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "'0.0983"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").Value = "'19.49"
$ws.Range("E47").Value = "  +1.52%  "
$ws.Range("D48").Value = "'0.0536"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "'10.27"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'0.0229"
$ws.Range("E50").Value = "  -1.16%  "
$ws.Range("D51").Value = "1.973.63"
$ws.Range("E51").Value = "  +1.02%  "
